# Auto-generated COM-interop script applying the 2023-09-14 data refresh
# to the violent-crime-full-year workbook. Updates partial-year 2023 (column J)
# figures plus a handful of small historical-year corrections, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5372
$ws.Range("J3").Value = 5697
$ws.Range("C4").Value = 1835
$ws.Range("E4").Value = 2006
$ws.Range("F4").Value = 1899
$ws.Range("J5").Value = 444
$ws.Range("J6").Value = 7124
$ws.Range("C7").Value = 28379
$ws.Range("E7").Value = 26011
$ws.Range("F7").Value = 24090
$ws.Range("J7").Value = 19891

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J4").Value = 9

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 346
$ws.Range("J3").Value = 381
$ws.Range("J4").Value = 73
$ws.Range("J6").Value = 418
$ws.Range("J7").Value = 1253

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 123
$ws.Range("J7").Value = 413

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 225
$ws.Range("J3").Value = 301
$ws.Range("J6").Value = 314
$ws.Range("J7").Value = 918

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 104
$ws.Range("J3").Value = 101
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 149
$ws.Range("J7").Value = 510

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 88
$ws.Range("J7").Value = 310

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 78
$ws.Range("J7").Value = 579
$ws.Range("J8").Value = 1253
$ws.Range("J11").Value = 309
$ws.Range("J18").Value = 167
$ws.Range("J19").Value = 574
$ws.Range("J20").Value = 415
$ws.Range("J23").Value = 191
$ws.Range("J25").Value = 100
$ws.Range("J29").Value = 1111
$ws.Range("J33").Value = 918
$ws.Range("J36").Value = 273
$ws.Range("J41").Value = 125
$ws.Range("J42").Value = 811
$ws.Range("J44").Value = 150
$ws.Range("J47").Value = 152
$ws.Range("J49").Value = 133
$ws.Range("J52").Value = 501
$ws.Range("J53").Value = 269
$ws.Range("C63").Value = 266
$ws.Range("E63").Value = 349
$ws.Range("F63").Value = 187
$ws.Range("J63").Value = 72
$ws.Range("J64").Value = 134
$ws.Range("J65").Value = 510
$ws.Range("J67").Value = 763
$ws.Range("J73").Value = 184
$ws.Range("J78").Value = 247
$ws.Range("J79").Value = 569
$ws.Range("J83").Value = 413
$ws.Range("J84").Value = 172
$ws.Range("J85").Value = 847
$ws.Range("J87").Value = 69
$ws.Range("J88").Value = 217
$ws.Range("J89").Value = 258
$ws.Range("J90").Value = 218
$ws.Range("J92").Value = 59
$ws.Range("J94").Value = 198
$ws.Range("J95").Value = 296
$ws.Range("J97").Value = 163
$ws.Range("J98").Value = 138
$ws.Range("J99").Value = 310
$ws.Range("J100").Value = 41
$ws.Range("C101").Value = 28379
$ws.Range("E101").Value = 26011
$ws.Range("F101").Value = 24090
$ws.Range("J101").Value = 19891

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 192
$ws.Range("J3").Value = 292
$ws.Range("J7").Value = 763

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 332
$ws.Range("J3").Value = 383
$ws.Range("J6").Value = 292
$ws.Range("J7").Value = 1111

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 143
$ws.Range("J3").Value = 170
$ws.Range("J6").Value = 211
$ws.Range("J7").Value = 574

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 125

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 183
$ws.Range("J3").Value = 163
$ws.Range("J6").Value = 410
$ws.Range("J7").Value = 811

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 81
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 62
$ws.Range("J4").Value = 9

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 161
$ws.Range("J3").Value = 201
$ws.Range("J7").Value = 569

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 115
$ws.Range("J4").Value = 39
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 415

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 273

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 181
$ws.Range("J7").Value = 579

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 44
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 40
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 138

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 60
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 83
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 258

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 81
$ws.Range("J5").Value = 7
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 220
$ws.Range("J3").Value = 312
$ws.Range("J4").Value = 57
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 847

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 149
$ws.Range("J6").Value = 204
$ws.Range("J7").Value = 501

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 69
